# Add a new "Driving Distance (miles)" column (F) to the "Driving Times"
# sheet, matching the header formatting already used by column E, and
# update row 2 with the new driving-time / driving-distance values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driving Times")

# New header cell F1, formatted like E1 (bold, centered, thin border).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Driving Distance (miles)"

# Updated driving time, and new driving distance, for row 2.
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 6.6524
